# Add the new "第五周 周二" (2017.9.26) plan block to the bottom of the sheet,
# mirroring the layout/formatting of the preceding weekly block (rows 82-90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the cell formatting (styles) of the previous week's block (rows 82-90)
# down onto the new block (rows 92-100), so fonts/fills/borders/number
# formats/alignment all match without having to rebuild the style table by
# hand.
$ws.Range("A82:D90").Copy()
$ws.Range("A92:D100").PasteSpecial(-4122)

# Header row for the new date block.
$ws.Range("A92").Value = "日期：2017.9.26 第五周 周二"

# Column headers (same as every other week's block).
$ws.Range("A93").Value = "人员"
$ws.Range("B93").Value = "计划任务"
$ws.Range("C93").Value = "完成情况"
$ws.Range("D93").Value = "备注"

# Per-person rows.
$ws.Range("A94").Value = "李杰"
$ws.Range("B94").Value = "熟悉小马后台框架并掌握操作数据库的方法"

$ws.Range("A95").Value = "周振朋"
$ws.Range("B95").Value = "继续编写“个人信息管理”用例规约"

$ws.Range("A96").Value = "禤锦辉"
$ws.Range("B96").Value = "继续编写“首页”用例规约"

$ws.Range("A97").Value = "柯新钿"
$ws.Range("B97").Value = "继续编写“账号管理”用例规约"

$ws.Range("A98").Value = "冯文雄"
$ws.Range("B98").Value = "熟悉小马后台框架并掌握操作数据库的方法"

$ws.Range("A99").Value = "阿卜力孜"
$ws.Range("B99").Value = "继续编写“查看买卖信息”用例规约"

# Summary row.
$ws.Range("A100").Value = "总结："

# Merge the header/summary rows across all four columns, as with every
# other block in the sheet.
$ws.Range("A92:D92").Merge()
$ws.Range("A100:D100").Merge()

# Move the selection/scroll position to the newly added block, like the
# author would have after typing the last cell.
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B94").Select()

# Restore the (approximate) saved window size recorded by the desktop
# Excel client that produced this edit.
$excel.ActiveWindow.Width = 28695
$excel.ActiveWindow.Height = 13050
